# Insert a new weekly data row for "Feria Lagunitas de Puerto Montt - Ciboulette"
# at row 172, pushing the existing rows 172:264 down to 173:265.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 172; this shifts rows 172-264 -> 173-265
# and copies formatting (incl. number format) from the row above.
$ws.Rows.Item(172).Insert()

# Populate the newly inserted row 172 with the new record's data.
$ws.Range("A172").Value = 4
$ws.Range("B172").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C172").Value = "Los Lagos"
$ws.Range("D172").Value = 44830
$ws.Range("D172").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E172").Value = 10
$ws.Range("F172").Value = 100112039
$ws.Range("G172").Value = "Ciboulette"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 80
$ws.Range("K172").Value = 3000
$ws.Range("L172").Value = 3000
$ws.Range("M172").Value = 3000
$ws.Range("N172").Value = "`$/docena de atados"
$ws.Range("O172").Value = "Región Metropolitana"
$ws.Range("P172").Value = 1000
$ws.Range("Q172").Value = 3
$ws.Range("R172").Value = "Hortaliza"
